$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "x" column (A) checkbox-style entries for rows 5, 6, 9, 10, 11
$ws.Range("A5").Value = "x"
$ws.Range("A6").Value = "x"
$ws.Range("A9").Value = "x"
$ws.Range("A10").Value = "x"

# A11 is a brand new cell in a row that previously had nothing in column A,
# so give it the same formatting (vertical-top alignment) as the rest of
# column A before setting its value.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A11").Value = "x"
